$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Sembol"
$ws.Cells.Item(2, 1).Value = "NVDA"
$ws.Cells.Item(3, 1).Value = "GOOGL"
$ws.Cells.Item(4, 1).Value = "GOOG"
$ws.Cells.Item(5, 1).Value = "MSFT"
$ws.Cells.Item(6, 1).Value = "AMZN"
$ws.Cells.Item(7, 1).Value = "META"
$ws.Cells.Item(8, 1).Value = "AVGO"
$ws.Cells.Item(9, 1).Value = "TSLA"
$ws.Cells.Item(10, 1).Value = "V"
$ws.Cells.Item(11, 1).Value = "ORCL"
$ws.Cells.Item(12, 1).Value = "BABA"
$ws.Cells.Item(13, 1).Value = "AMD"
$ws.Cells.Item(14, 1).Value = "PLTR"
$ws.Cells.Item(15, 1).Value = "NFLX"
$ws.Cells.Item(16, 1).Value = "GS"
$ws.Cells.Item(17, 1).Value = "JPM-PD"
$ws.Cells.Item(18, 1).Value = "JPM-PC"
$ws.Cells.Item(19, 1).Value = "UNH"
$ws.Cells.Item(20, 1).Value = "MS"
$ws.Cells.Item(21, 1).Value = "BAC-PK"
$ws.Cells.Item(22, 1).Value = "BAC-PL"
$ws.Cells.Item(23, 1).Value = "NVO"
$ws.Cells.Item(24, 1).Value = "SAP"
$ws.Cells.Item(25, 1).Value = "IBM"
$ws.Cells.Item(26, 1).Value = "BAC-PE"
$ws.Cells.Item(27, 1).Value = "BML-PL"
$ws.Cells.Item(28, 1).Value = "BAC-PB"
$ws.Cells.Item(29, 1).Value = "AXP"
$ws.Cells.Item(30, 1).Value = "TMO"
$ws.Cells.Item(31, 1).Value = "CRM"
$ws.Cells.Item(32, 1).Value = "KLAC"
$ws.Cells.Item(33, 1).Value = "WFC-PY"
$ws.Cells.Item(34, 1).Value = "DIS"
$ws.Cells.Item(35, 1).Value = "WFC-PL"
$ws.Cells.Item(36, 1).Value = "APH"
$ws.Cells.Item(37, 1).Value = "ISRG"
$ws.Cells.Item(38, 1).Value = "ABT"
$ws.Cells.Item(39, 1).Value = "BX"
$ws.Cells.Item(40, 1).Value = "APP"
$ws.Cells.Item(41, 1).Value = "ANET"
$ws.Cells.Item(42, 1).Value = "SHOP"
$ws.Cells.Item(43, 1).Value = "ACN"
$ws.Cells.Item(44, 1).Value = "BLK"
$ws.Cells.Item(45, 1).Value = "UBER"
$ws.Cells.Item(46, 1).Value = "DHR"
$ws.Cells.Item(47, 1).Value = "BKNG"
$ws.Cells.Item(48, 1).Value = "QCOM"
$ws.Cells.Item(49, 1).Value = "SPGI"
$ws.Cells.Item(50, 1).Value = "INTU"
$ws.Cells.Item(51, 1).Value = "UBS"
$ws.Cells.Item(52, 1).Value = "PDD"
$ws.Cells.Item(53, 1).Value = "BBVA"
$ws.Cells.Item(54, 1).Value = "NOW"
$ws.Cells.Item(55, 1).Value = "COF"
$ws.Cells.Item(56, 1).Value = "BSX"
$ws.Cells.Item(57, 1).Value = "NEM"
$ws.Cells.Item(58, 1).Value = "SONY"
$ws.Cells.Item(59, 1).Value = "SYK"
$ws.Cells.Item(60, 1).Value = "PANW"
$ws.Cells.Item(61, 1).Value = "IBKR"
$ws.Cells.Item(62, 1).Value = "ADBE"
$ws.Cells.Item(63, 1).Value = "WFC-PC"
$ws.Cells.Item(64, 1).Value = "PGR"
$ws.Cells.Item(65, 1).Value = "CRWD"
$ws.Cells.Item(66, 1).Value = "MELI"
$ws.Cells.Item(67, 1).Value = "KKR"
$ws.Cells.Item(68, 1).Value = "AEM"
$ws.Cells.Item(69, 1).Value = "BN"
$ws.Cells.Item(70, 1).Value = "SPOT"
$ws.Cells.Item(71, 1).Value = "ADP"
$ws.Cells.Item(72, 1).Value = "CEG"
$ws.Cells.Item(73, 1).Value = "CVNA"
$ws.Cells.Item(74, 1).Value = "SNPS"
$ws.Cells.Item(75, 1).Value = "HOOD"
$ws.Cells.Item(76, 1).Value = "MCO"
$ws.Cells.Item(77, 1).Value = "DASH"
$ws.Cells.Item(78, 1).Value = "NU"
$ws.Cells.Item(79, 1).Value = "CDNS"
$ws.Cells.Item(80, 1).Value = "B"
$ws.Cells.Item(81, 1).Value = "NTES"
$ws.Cells.Item(82, 1).Value = "ELV"
$ws.Cells.Item(83, 1).Value = "ORLY"
$ws.Cells.Item(84, 1).Value = "MS-PK"
$ws.Cells.Item(85, 1).Value = "AMT"
$ws.Cells.Item(86, 1).Value = "BAM"
$ws.Cells.Item(87, 1).Value = "MS-PI"
$ws.Cells.Item(88, 1).Value = "ABNB"
$ws.Cells.Item(89, 1).Value = "TDG"
$ws.Cells.Item(90, 1).Value = "DB"
$ws.Cells.Item(91, 1).Value = "MS-PF"
$ws.Cells.Item(92, 1).Value = "MS-PE"
$ws.Cells.Item(93, 1).Value = "CMI"
$ws.Cells.Item(94, 1).Value = "APO"
$ws.Cells.Item(95, 1).Value = "INFY"
$ws.Cells.Item(96, 1).Value = "SE"
$ws.Cells.Item(97, 1).Value = "USB-PP"
$ws.Cells.Item(98, 1).Value = "AON"
$ws.Cells.Item(99, 1).Value = "SNOW"
$ws.Cells.Item(100, 1).Value = "RELX"
$ws.Cells.Item(101, 1).Value = "MRVL"
$ws.Cells.Item(102, 1).Value = "WBD"
$ws.Cells.Item(103, 1).Value = "TEL"
$ws.Cells.Item(104, 1).Value = "NET"
$ws.Cells.Item(105, 1).Value = "AJG"
$ws.Cells.Item(106, 1).Value = "AZO"
$ws.Cells.Item(107, 1).Value = "DUK-PA"
$ws.Cells.Item(108, 1).Value = "RACE"
$ws.Cells.Item(109, 1).Value = "RKT"
$ws.Cells.Item(110, 1).Value = "CTA-PB"
$ws.Cells.Item(111, 1).Value = "NXPI"
$ws.Cells.Item(112, 1).Value = "ADSK"
$ws.Cells.Item(113, 1).Value = "COIN"
$ws.Cells.Item(114, 1).Value = "NDAQ"
$ws.Cells.Item(115, 1).Value = "SRE"
$ws.Cells.Item(116, 1).Value = "IDXX"
$ws.Cells.Item(117, 1).Value = "TRI"
$ws.Cells.Item(118, 1).Value = "BIDU"
$ws.Cells.Item(119, 1).Value = "PYPL"
$ws.Cells.Item(120, 1).Value = "VST"
$ws.Cells.Item(121, 1).Value = "F"
$ws.Cells.Item(122, 1).Value = "CCJ"
$ws.Cells.Item(123, 1).Value = "RBLX"
$ws.Cells.Item(124, 1).Value = "ARGX"
$ws.Cells.Item(125, 1).Value = "MET"
$ws.Cells.Item(126, 1).Value = "EA"
$ws.Cells.Item(127, 1).Value = "SCHW-PD"
$ws.Cells.Item(128, 1).Value = "WDAY"
$ws.Cells.Item(129, 1).Value = "BSBR"
$ws.Cells.Item(130, 1).Value = "FNV"
$ws.Cells.Item(131, 1).Value = "EW"
$ws.Cells.Item(132, 1).Value = "ARES"
$ws.Cells.Item(133, 1).Value = "CRWV"
$ws.Cells.Item(134, 1).Value = "AXON"
$ws.Cells.Item(135, 1).Value = "DDOG"
$ws.Cells.Item(136, 1).Value = "ALNY"
$ws.Cells.Item(137, 1).Value = "ROK"
$ws.Cells.Item(138, 1).Value = "CTA-PA"
$ws.Cells.Item(139, 1).Value = "HEI"
$ws.Cells.Item(140, 1).Value = "MSTR"
$ws.Cells.Item(141, 1).Value = "KGC"
$ws.Cells.Item(142, 1).Value = "MSCI"
$ws.Cells.Item(143, 1).Value = "TTWO"
$ws.Cells.Item(144, 1).Value = "SPG-PJ"
$ws.Cells.Item(145, 1).Value = "EXC"
$ws.Cells.Item(146, 1).Value = "XEL"
$ws.Cells.Item(147, 1).Value = "ROP"
$ws.Cells.Item(148, 1).Value = "JD"
$ws.Cells.Item(149, 1).Value = "EBAY"
$ws.Cells.Item(150, 1).Value = "MET-PE"
$ws.Cells.Item(151, 1).Value = "RKLB"
$ws.Cells.Item(152, 1).Value = "MET-PA"
$ws.Cells.Item(153, 1).Value = "EL"
$ws.Cells.Item(154, 1).Value = "CTSH"
$ws.Cells.Item(155, 1).Value = "TCOM"
$ws.Cells.Item(156, 1).Value = "LVS"
$ws.Cells.Item(157, 1).Value = "IQV"
$ws.Cells.Item(158, 1).Value = "RDDT"
$ws.Cells.Item(159, 1).Value = "MCHP"
$ws.Cells.Item(160, 1).Value = "XYZ"
$ws.Cells.Item(161, 1).Value = "HEI-A"
$ws.Cells.Item(162, 1).Value = "ALC"
$ws.Cells.Item(163, 1).Value = "GRMN"
$ws.Cells.Item(164, 1).Value = "ASTS"
$ws.Cells.Item(165, 1).Value = "A"
$ws.Cells.Item(166, 1).Value = "PRU"
$ws.Cells.Item(167, 1).Value = "PSA-PH"
$ws.Cells.Item(168, 1).Value = "PAYX"
$ws.Cells.Item(169, 1).Value = "CCI"
$ws.Cells.Item(170, 1).Value = "MDLN"
$ws.Cells.Item(171, 1).Value = "FICO"
$ws.Cells.Item(172, 1).Value = "VEEV"
$ws.Cells.Item(173, 1).Value = "FISV"
$ws.Cells.Item(174, 1).Value = "TEAM"
$ws.Cells.Item(175, 1).Value = "RYAAY"
$ws.Cells.Item(176, 1).Value = "CPNG"
$ws.Cells.Item(177, 1).Value = "SYM"
$ws.Cells.Item(178, 1).Value = "APO-PA"
$ws.Cells.Item(179, 1).Value = "SATS"
$ws.Cells.Item(180, 1).Value = "LYV"
$ws.Cells.Item(181, 1).Value = "ZS"
$ws.Cells.Item(182, 1).Value = "INSM"
$ws.Cells.Item(183, 1).Value = "NTRA"
$ws.Cells.Item(184, 1).Value = "MDB"
$ws.Cells.Item(185, 1).Value = "EXPE"
$ws.Cells.Item(186, 1).Value = "CHT"
$ws.Cells.Item(187, 1).Value = "UI"
$ws.Cells.Item(188, 1).Value = "BE"
$ws.Cells.Item(189, 1).Value = "ESLT"
$ws.Cells.Item(190, 1).Value = "PSA-PK"
$ws.Cells.Item(191, 1).Value = "ALL-PH"
$ws.Cells.Item(192, 1).Value = "SOFI"
$ws.Cells.Item(193, 1).Value = "ALL-PB"
$ws.Cells.Item(194, 1).Value = "FOXA"
$ws.Cells.Item(195, 1).Value = "HUM"
$ws.Cells.Item(196, 1).Value = "EXR"
$ws.Cells.Item(197, 1).Value = "FIS"
$ws.Cells.Item(198, 1).Value = "FOX"
$ws.Cells.Item(199, 1).Value = "VRSK"
$ws.Cells.Item(200, 1).Value = "FLUT"
$ws.Cells.Item(201, 1).Value = "BNTX"
$ws.Cells.Item(202, 1).Value = "MTD"
$ws.Cells.Item(203, 1).Value = "NRG"
$ws.Cells.Item(204, 1).Value = "SYF"
$ws.Cells.Item(205, 1).Value = "DXCM"
$ws.Cells.Item(206, 1).Value = "TME"
$ws.Cells.Item(207, 1).Value = "STLA"
$ws.Cells.Item(208, 1).Value = "CSGP"
$ws.Cells.Item(209, 1).Value = "ALAB"
$ws.Cells.Item(210, 1).Value = "PAAS"
$ws.Cells.Item(211, 1).Value = "WIT"
$ws.Cells.Item(212, 1).Value = "BRO"
$ws.Cells.Item(213, 1).Value = "EFX"
$ws.Cells.Item(214, 1).Value = "ES"
$ws.Cells.Item(215, 1).Value = "FSLR"
$ws.Cells.Item(216, 1).Value = "STE"
$ws.Cells.Item(217, 1).Value = "AER"
$ws.Cells.Item(218, 1).Value = "DLTR"
$ws.Cells.Item(219, 1).Value = "AWK"
$ws.Cells.Item(220, 1).Value = "OMC"
$ws.Cells.Item(221, 1).Value = "AVB"
$ws.Cells.Item(222, 1).Value = "VLTO"
$ws.Cells.Item(223, 1).Value = "DLR-PK"
$ws.Cells.Item(224, 1).Value = "RGLD"
$ws.Cells.Item(225, 1).Value = "BR"
$ws.Cells.Item(226, 1).Value = "PSTG"
$ws.Cells.Item(227, 1).Value = "SQM"
$ws.Cells.Item(228, 1).Value = "FLEX"
$ws.Cells.Item(229, 1).Value = "AXIA-PC"
$ws.Cells.Item(230, 1).Value = "ILMN"
$ws.Cells.Item(231, 1).Value = "VRSN"
$ws.Cells.Item(232, 1).Value = "TPG"
$ws.Cells.Item(233, 1).Value = "TROW"
$ws.Cells.Item(234, 1).Value = "WAT"
$ws.Cells.Item(235, 1).Value = "CRDO"
$ws.Cells.Item(236, 1).Value = "NBIS"
$ws.Cells.Item(237, 1).Value = "LULU"
$ws.Cells.Item(238, 1).Value = "OWL"
$ws.Cells.Item(239, 1).Value = "CNC"
$ws.Cells.Item(240, 1).Value = "FUTU"
$ws.Cells.Item(241, 1).Value = "AFRM"
$ws.Cells.Item(242, 1).Value = "DLR-PJ"
$ws.Cells.Item(243, 1).Value = "TLK"
$ws.Cells.Item(244, 1).Value = "CYBR"
$ws.Cells.Item(245, 1).Value = "FWONK"
$ws.Cells.Item(246, 1).Value = "ALB"
$ws.Cells.Item(247, 1).Value = "PSLV"
$ws.Cells.Item(248, 1).Value = "FWONA"
$ws.Cells.Item(249, 1).Value = "CG"
$ws.Cells.Item(250, 1).Value = "RL"
$ws.Cells.Item(251, 1).Value = "GPN"
$ws.Cells.Item(252, 1).Value = "AS"
$ws.Cells.Item(253, 1).Value = "SSNC"
$ws.Cells.Item(254, 1).Value = "GMAB"
$ws.Cells.Item(255, 1).Value = "TWLO"
$ws.Cells.Item(256, 1).Value = "Q"
$ws.Cells.Item(257, 1).Value = "HL"
$ws.Cells.Item(258, 1).Value = "SBAC"
$ws.Cells.Item(259, 1).Value = "RCI"
$ws.Cells.Item(260, 1).Value = "CHKP"
$ws.Cells.Item(261, 1).Value = "PTC"
$ws.Cells.Item(262, 1).Value = "TOST"
$ws.Cells.Item(263, 1).Value = "GIB"
$ws.Cells.Item(264, 1).Value = "RIVN"
$ws.Cells.Item(265, 1).Value = "PODD"
$ws.Cells.Item(266, 1).Value = "TYL"
$ws.Cells.Item(267, 1).Value = "RVMD"
$ws.Cells.Item(268, 1).Value = "BWXT"
$ws.Cells.Item(269, 1).Value = "KTOS"
$ws.Cells.Item(270, 1).Value = "MRNA"
$ws.Cells.Item(271, 1).Value = "HIG-PG"
$ws.Cells.Item(272, 1).Value = "GRAB"
$ws.Cells.Item(273, 1).Value = "IOT"
$ws.Cells.Item(274, 1).Value = "CX"
$ws.Cells.Item(275, 1).Value = "U"
$ws.Cells.Item(276, 1).Value = "HPQ"
$ws.Cells.Item(277, 1).Value = "CRCL"
$ws.Cells.Item(278, 1).Value = "FITBI"
$ws.Cells.Item(279, 1).Value = "XPEV"
$ws.Cells.Item(280, 1).Value = "IT"
$ws.Cells.Item(281, 1).Value = "PSNYW"
$ws.Cells.Item(282, 1).Value = "AGI"
$ws.Cells.Item(283, 1).Value = "NVT"
$ws.Cells.Item(284, 1).Value = "ALLY"
$ws.Cells.Item(285, 1).Value = "PNR"
$ws.Cells.Item(286, 1).Value = "PINS"
$ws.Cells.Item(287, 1).Value = "SN"
$ws.Cells.Item(288, 1).Value = "WST"
$ws.Cells.Item(289, 1).Value = "HUBS"
$ws.Cells.Item(290, 1).Value = "NWS"
$ws.Cells.Item(291, 1).Value = "IREN"
$ws.Cells.Item(292, 1).Value = "ZG"
$ws.Cells.Item(293, 1).Value = "TRMB"
$ws.Cells.Item(294, 1).Value = "MEDP"
$ws.Cells.Item(295, 1).Value = "JLL"
$ws.Cells.Item(296, 1).Value = "QXO"
$ws.Cells.Item(297, 1).Value = "Z"
$ws.Cells.Item(298, 1).Value = "TRU"
$ws.Cells.Item(299, 1).Value = "TTD"
$ws.Cells.Item(300, 1).Value = "INVH"
$ws.Cells.Item(301, 1).Value = "CDE"
$ws.Cells.Item(302, 1).Value = "NLY"
$ws.Cells.Item(303, 1).Value = "HII"
$ws.Cells.Item(304, 1).Value = "MAA"
$ws.Cells.Item(305, 1).Value = "OKTA"
$ws.Cells.Item(306, 1).Value = "SUI"
$ws.Cells.Item(307, 1).Value = "TLN"
$ws.Cells.Item(308, 1).Value = "GFL"
$ws.Cells.Item(309, 1).Value = "KSPI"
$ws.Cells.Item(310, 1).Value = "ROKU"
$ws.Cells.Item(311, 1).Value = "GEN"
$ws.Cells.Item(312, 1).Value = "NWSA"
$ws.Cells.Item(313, 1).Value = "IONQ"
$ws.Cells.Item(314, 1).Value = "ONON"
$ws.Cells.Item(315, 1).Value = "DKNG"
$ws.Cells.Item(316, 1).Value = "AA"
$ws.Cells.Item(317, 1).Value = "AVAV"
$ws.Cells.Item(318, 1).Value = "KEY-PI"
$ws.Cells.Item(319, 1).Value = "BBIO"
$ws.Cells.Item(320, 1).Value = "GH"
$ws.Cells.Item(321, 1).Value = "MLI"
$ws.Cells.Item(322, 1).Value = "HMY"
$ws.Cells.Item(323, 1).Value = "PFGC"
$ws.Cells.Item(324, 1).Value = "ULS"
$ws.Cells.Item(325, 1).Value = "GDDY"
$ws.Cells.Item(326, 1).Value = "ARCC"
$ws.Cells.Item(327, 1).Value = "ICLR"
$ws.Cells.Item(328, 1).Value = "FIG"
$ws.Cells.Item(329, 1).Value = "RGC"
$ws.Cells.Item(330, 1).Value = "KRMN"
$ws.Cells.Item(331, 1).Value = "RF-PC"
$ws.Cells.Item(332, 1).Value = "W"
$ws.Cells.Item(333, 1).Value = "ASND"
$ws.Cells.Item(334, 1).Value = "CACI"
$ws.Cells.Item(335, 1).Value = "AKAM"
$ws.Cells.Item(336, 1).Value = "JHX"
$ws.Cells.Item(337, 1).Value = "PEN"
$ws.Cells.Item(338, 1).Value = "EVR"
$ws.Cells.Item(339, 1).Value = "CELH"
$ws.Cells.Item(340, 1).Value = "DPZ"
$ws.Cells.Item(341, 1).Value = "EMBJ"
$ws.Cells.Item(342, 1).Value = "SBSW"
$ws.Cells.Item(343, 1).Value = "LOGI"
$ws.Cells.Item(344, 1).Value = "EQH"
$ws.Cells.Item(345, 1).Value = "GWRE"
$ws.Cells.Item(346, 1).Value = "SOLV"
$ws.Cells.Item(347, 1).Value = "BILI"
$ws.Cells.Item(348, 1).Value = "RBRK"
$ws.Cells.Item(349, 1).Value = "FIGR"
$ws.Cells.Item(350, 1).Value = "AMH"
$ws.Cells.Item(351, 1).Value = "RVTY"
$ws.Cells.Item(352, 1).Value = "JKHY"
$ws.Cells.Item(353, 1).Value = "RYAN"
$ws.Cells.Item(354, 1).Value = "PSKY"
$ws.Cells.Item(355, 1).Value = "CHWY"
$ws.Cells.Item(356, 1).Value = "UNM"
$ws.Cells.Item(357, 1).Value = "SNAP"
$ws.Cells.Item(358, 1).Value = "JEF"
$ws.Cells.Item(359, 1).Value = "BNT"
$ws.Cells.Item(360, 1).Value = "OKLO"
$ws.Cells.Item(361, 1).Value = "HLI"
$ws.Cells.Item(362, 1).Value = "EQX"
$ws.Cells.Item(363, 1).Value = "IVZ"
$ws.Cells.Item(364, 1).Value = "AGNC"
$ws.Cells.Item(365, 1).Value = "BMNR"
$ws.Cells.Item(366, 1).Value = "GLXY"
$ws.Cells.Item(367, 1).Value = "AG"
$ws.Cells.Item(368, 1).Value = "GMED"
$ws.Cells.Item(369, 1).Value = "AMKR"
$ws.Cells.Item(370, 1).Value = "DT"
$ws.Cells.Item(371, 1).Value = "ACGLO"
$ws.Cells.Item(372, 1).Value = "DOC"
$ws.Cells.Item(373, 1).Value = "SMMT"
$ws.Cells.Item(374, 1).Value = "RMBS"
$ws.Cells.Item(375, 1).Value = "EPAM"
$ws.Cells.Item(376, 1).Value = "JOBY"
$ws.Cells.Item(377, 1).Value = "CMA"
$ws.Cells.Item(378, 1).Value = "TEM"
$ws.Cells.Item(379, 1).Value = "NYT"
$ws.Cells.Item(380, 1).Value = "NTNX"
$ws.Cells.Item(381, 1).Value = "BSY"
$ws.Cells.Item(382, 1).Value = "DOCU"
$ws.Cells.Item(383, 1).Value = "CPT"
$ws.Cells.Item(384, 1).Value = "BXP"
$ws.Cells.Item(385, 1).Value = "CEF"
$ws.Cells.Item(386, 1).Value = "MDGL"
$ws.Cells.Item(387, 1).Value = "BAH"
$ws.Cells.Item(388, 1).Value = "MICC"
$ws.Cells.Item(389, 1).Value = "MP"
$ws.Cells.Item(390, 1).Value = "QGEN"
$ws.Cells.Item(391, 1).Value = "WTRG"
$ws.Cells.Item(392, 1).Value = "SARO"
$ws.Cells.Item(393, 1).Value = "CRL"
$ws.Cells.Item(394, 1).Value = "UHAL"
$ws.Cells.Item(395, 1).Value = "MOH"
$ws.Cells.Item(396, 1).Value = "DRS"
$ws.Cells.Item(397, 1).Value = "MANH"
$ws.Cells.Item(398, 1).Value = "FDS"
$ws.Cells.Item(399, 1).Value = "AFG"
$ws.Cells.Item(400, 1).Value = "CART"
$ws.Cells.Item(401, 1).Value = "SEIC"
$ws.Cells.Item(402, 1).Value = "CAE"
$ws.Cells.Item(403, 1).Value = "YMM"
$ws.Cells.Item(404, 1).Value = "XP"
$ws.Cells.Item(405, 1).Value = "PCOR"
$ws.Cells.Item(406, 1).Value = "KLAR"
$ws.Cells.Item(407, 1).Value = "APLD"
$ws.Cells.Item(408, 1).Value = "JAZZ"
$ws.Cells.Item(409, 1).Value = "CHYM"
$ws.Cells.Item(410, 1).Value = "VNO-PL"
$ws.Cells.Item(411, 1).Value = "UHAL-B"
$ws.Cells.Item(412, 1).Value = "SAIL"
$ws.Cells.Item(413, 1).Value = "SANM"
$ws.Cells.Item(414, 1).Value = "NGD"
$ws.Cells.Item(415, 1).Value = "VNO-PM"
$ws.Cells.Item(416, 1).Value = "COMP"
$ws.Cells.Item(417, 1).Value = "BIO-B"
$ws.Cells.Item(418, 1).Value = "SOLS"
$ws.Cells.Item(419, 1).Value = "ARE"
$ws.Cells.Item(420, 1).Value = "EGO"
$ws.Cells.Item(421, 1).Value = "BROS"
$ws.Cells.Item(422, 1).Value = "REXR"
$ws.Cells.Item(423, 1).Value = "ABVX"
$ws.Cells.Item(424, 1).Value = "UWMC"
$ws.Cells.Item(425, 1).Value = "ARWR"
$ws.Cells.Item(426, 1).Value = "AXSM"
$ws.Cells.Item(427, 1).Value = "RGEN"
$ws.Cells.Item(428, 1).Value = "DOX"
$ws.Cells.Item(429, 1).Value = "STEP"
$ws.Cells.Item(430, 1).Value = "MORN"
$ws.Cells.Item(431, 1).Value = "AMG"
$ws.Cells.Item(432, 1).Value = "UEC"
$ws.Cells.Item(433, 1).Value = "LUMN"
$ws.Cells.Item(434, 1).Value = "GGAL"
$ws.Cells.Item(435, 1).Value = "QBTS"
$ws.Cells.Item(436, 1).Value = "TTAN"
$ws.Cells.Item(437, 1).Value = "RZB"
$ws.Cells.Item(438, 1).Value = "AGNCM"
$ws.Cells.Item(439, 1).Value = "AGNCN"
$ws.Cells.Item(440, 1).Value = "PEGA"
$ws.Cells.Item(441, 1).Value = "OR"
$ws.Cells.Item(442, 1).Value = "AUR"
$ws.Cells.Item(443, 1).Value = "GGB"
$ws.Cells.Item(444, 1).Value = "UGI"
$ws.Cells.Item(445, 1).Value = "PL"
$ws.Cells.Item(446, 1).Value = "LEVI"

$ws.Range("A447:A529").ClearContents()
